$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 <- original row 24 content
$ws.Range("A23").Value = 102941038
$ws.Range("B23").Value = 56395
$ws.Range("C23").Value = "Ovaliderad"
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = "Tretåig hackspett"
$ws.Range("G23").Value = "Picoides tridactylus"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
$ws.Range("I23").Value = "1"
$ws.Range("K23").Value = ""
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = "äldre spår"
$ws.Range("P23").Value = "Gyljberget, nord ost, Hls"
$ws.Range("Q23").Value = 551737.9524457334
$ws.Range("R23").Value = 6765756.293540224
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = "Gävleborg"
$ws.Range("U23").Value = "Bollnäs"
$ws.Range("V23").Value = "Hälsingland"
$ws.Range("W23").Value = "Bollnäs"
$ws.Range("Y23").Value = "2022-08-17"
$ws.Range("Z23").Value = "10:52"
$ws.Range("AA23").Value = "2022-08-17"
$ws.Range("AB23").Value = "10:52"
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AT23").Value = ""
$ws.Range("AW23").Value = "Enar Gesslin"
$ws.Range("AX23").Value = "Enar Gesslin"
$ws.Range("AY23").Value = ""

# Row 24 <- original row 23 content
$ws.Range("A24").Value = 97604899
$ws.Range("B24").Value = 77506
$ws.Range("C24").Value = "Ovaliderad"
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 6425
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("I24").Value = ""
$ws.Range("P24").Value = "Gyljberget–Igeltjärnen, Hls"
$ws.Range("Q24").Value = 551888.2619905178
$ws.Range("R24").Value = 6765241.166833818
$ws.Range("S24").Value = 25
$ws.Range("T24").Value = "Gävleborg"
$ws.Range("U24").Value = "Bollnäs"
$ws.Range("V24").Value = "Hälsingland"
$ws.Range("W24").Value = "Bollnäs"
$ws.Range("Y24").Value = "2021-05-18"
$ws.Range("Z24").Value = "00:00"
$ws.Range("AA24").Value = "2021-05-18"
$ws.Range("AB24").Value = "00:00"
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AT24").Value = ""
$ws.Range("AW24").Value = "Jenny Andersson"
$ws.Range("AX24").Value = "Jenny Andersson"
$ws.Range("AY24").Value = "Länsstyrelsen Gävleborg funktionsindelning"
